# feat: add struct to excel
#
# - Remove the now-unused "Sheet1" worksheet.
# - On "Infos", insert a new header/struct row (row 2) that duplicates the
#   existing header values from row 1 (user_name / phone / age / man),
#   pushing the two data rows down to rows 3-4.
# - Make "Infos" the active/selected sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets("Infos")

# Insert a new row right below the header row; existing rows 2-3 shift to 3-4.
$ws.Rows("2:2").Insert()

# Populate the newly inserted struct/header row with the same field names.
$ws.Range("A2").Value = "user_name"
$ws.Range("B2").Value = "phone"
$ws.Range("C2").Value = "age"
$ws.Range("D2").Value = "man"

# Drop the now-empty default sheet.
$null = $wb.Worksheets("Sheet1").Delete()

# Make "Infos" the selected/active tab.
$ws.Activate()
